$wb = $excel.ActiveWorkbook

# Rename sheets (tab names encode a task-order timestamp)
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477877952398"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778805413942"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778805423946"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477880590393"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778806533933"

# Sheet 1 - GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778779193935.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778779354267.csv"
$ws1.Range("B4").Value = "go_stims-16504778779363961.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477877952398.csv"

# Sheet 2 - NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_5-165047787805843.csv"
$ws2.Range("B3").Value = "ZB-match_1-1650477878005424.csv"
$ws2.Range("B4").Value = "OB-16504778782803905.csv"
$ws2.Range("B5").Value = "TB-16504778796104274.csv"
$ws2.Range("B6").Value = "OB-16504778785763981.csv"
$ws2.Range("B7").Value = "TB-16504778800234234.csv"
$ws2.Range("B8").Value = "ZB-match_6-16504778781894267.csv"
$ws2.Range("B9").Value = "OB-16504778784133923.csv"
$ws2.Range("B10").Value = "TB-16504778805184252.csv"

# Sheet 3 - RS (no content changes)

# Sheet 4 - TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778805573933.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778805433948.csv"
$ws4.Range("B4").Value = "MM_stims-16504778805734258.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778805573933.csv"
$ws4.Range("B6").Value = "MM_stims-16504778805894253.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778805734258.csv"

# Sheet 5 - vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650477880605426.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778805934374.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778806214263.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778806374266.csv"
